# "functioneaza cum trebuie citirea / scrierea goalurilor EXCEL"
#
# The sheet used to hold a generic two-goal template (Goal1/Task1.x,
# Goal2/Task2.x) followed by a couple of sample "Goal" rows and a merged
# "FINISH" banner row. The app now reads/writes real goal data, so the
# template rows are replaced by the actual exported goals/tasks and the
# trailing banner + spacer rows are dropped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Helper: write a literal text value (bypassing Excel's automatic
# number/percentage parsing) without leaving a lingering style behind.
# ---------------------------------------------------------------------
function Set-LiteralText {
    param($cell, [string]$text)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# ---------------------------------------------------------------------
# 1. Drop the old trailing rows: two blank spacer rows, the old sample
#    "Goal" rows and the merged "FINISH" footer banner (rows 14-18).
# ---------------------------------------------------------------------
$ws.Rows("14:18").Delete()

# Rows 3 & 4 used to carry a trailing "<computed>" marker in column G
# that doesn't exist for the new data - drop the stale leftovers.
$ws.Cells.Item(3, 7).Clear()
$ws.Cells.Item(4, 7).Clear()

# ---------------------------------------------------------------------
# 2. Headers (rows 1-2) - values are unchanged, only the merged-cell
#    style index shifts because of new styles inserted earlier in the
#    table; re-assert the values so everything round-trips cleanly.
# ---------------------------------------------------------------------
$ws.Range("B1").Value = "Description"
$ws.Range("D1").Value = "Estimated Time"
$ws.Range("E1").Value = "Realized Time "
$ws.Range("F1").Value = "Progress"

$ws.Range("B2").Value = "Description"
$ws.Range("C2").Value = "Estimated Date"
$ws.Range("D2").Value = "Estimated Time"
$ws.Range("E2").Value = "Realized Time "
$ws.Range("F2").Value = "Value"
$ws.Range("G2").Value = "Label"

# ---------------------------------------------------------------------
# 3. Data rows 3-13: the real goal + tasks, then every goal the user
#    actually created in the app.
# ---------------------------------------------------------------------
function Set-GoalRow {
    param($row, $a, $b, $c, $label)

    $ws.Cells.Item($row, 1).Value = $a
    $ws.Cells.Item($row, 2).Value = $b

    $ws.Cells.Item($row, 3).NumberFormat = "dd/mm/yyyy"
    $ws.Cells.Item($row, 3).Value = $c

    Set-LiteralText $ws.Cells.Item($row, 4) "0|0"
    Set-LiteralText $ws.Cells.Item($row, 5) "0|0"
    Set-LiteralText $ws.Cells.Item($row, 6) "0%"

    if ($label -ne $null) {
        $ws.Cells.Item($row, 7).Value = $label
    }
}

Set-GoalRow 3  "Goal2"   "I will make StandApp in 2 months"    44241.456400462965 $null
Set-GoalRow 4  "Task2.1" "do the view"                         44238.456400462965 $null
Set-GoalRow 5  "Task2.2" "do the model"                        44239.456400462965 $null
Set-GoalRow 6  "Task2.n" "finish this damn app"                44236.456396574074 $null
Set-GoalRow 7  "Goal"    "description"                         44236.456396574074 $null
Set-GoalRow 8  "Goal"    "description"                         44236.456396574074 $null
Set-GoalRow 9  "Goal"    "fafsafa"                              44237.647051111111 $null
Set-GoalRow 10 "Goal"    "abc"                                  44238.650255671295 $null
Set-GoalRow 11 "Goal"    "abc2"                                 44239.650395219905 $null
Set-GoalRow 12 "Goal"    "VREAU SA FIU ARTIST"                  44247.652703865744 $null

# Row 13 is the newest goal; its date cell uses the new dd/MM/yyyy
# format and it carries a "Label" (Just Started) in column G.
$ws.Cells.Item(13, 1).Value = "Goal"
$ws.Cells.Item(13, 2).Value = "brandNewGoal -> IULI CONDUCE LUMEEEEA"
$ws.Cells.Item(13, 3).NumberFormat = "dd/MM/yyyy"
$ws.Cells.Item(13, 3).Value = 44237.72078783565
Set-LiteralText $ws.Cells.Item(13, 4) "0|0"
Set-LiteralText $ws.Cells.Item(13, 5) "0|0"
Set-LiteralText $ws.Cells.Item(13, 6) "0%"
$ws.Cells.Item(13, 7).Value = "Just Started"

# ---------------------------------------------------------------------
# 4. A3 (the active goal row) gets a plain white solid fill highlight.
# ---------------------------------------------------------------------
$ws.Range("A3").Interior.ThemeColor = 2

# ---------------------------------------------------------------------
# 5. Selection matches the freshly-trimmed sheet (rows 13-17, which is
#    now past the used range).
# ---------------------------------------------------------------------
$ws.Range("A13:XFD17").Select()
